$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.424.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "'2.311.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'311.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'103.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.18%  "
$ws.Range("D7").Value = "'0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D10").Value = "'36.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("D11").Value = "'0.0816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").Value = "'51.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D15").Value = "'2.669.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "'15.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").Value = "'2.445.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.04%  "
$ws.Range("D18").Value = "'0.811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "'43.321.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").Value = "'12.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("D23").Value = "'68.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'242.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'24.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.17%  "
$ws.Range("E29").Value = "  +8.31%  "
$ws.Range("D30").Value = "'36.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'9.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").Value = "'168.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "'5.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  +6.21%  "
$ws.Range("D37").Value = "'0.0744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "'4.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.99%  "
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Value = "'2.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.73%  "
$ws.Range("E44").Value = "  +6.25%  "
$ws.Range("D45").Value = "'1.985.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'19.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "'3.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("D48").Value = "'9.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "'55.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.70%  "
$ws.Range("D50").Value = "'2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'1.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.31%  "
